# Apply the "Global M2 - Egypt" monthly data refresh.
# Updates a handful of restated historical values (rows 318-320, 327, 329)
# and appends a new observation row (330) for 2023-06-01.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revisions to existing rows ---

# Row 318 (2022-06-01)
$ws.Range("B318").Value = 6614488000000
$ws.Range("D318").Value = 352584648187.6332

# Row 319 (2022-07-01)
$ws.Range("B319").Value = 6708871000000
$ws.Range("D319").Value = 355531054583.9957

# Row 320 (2022-08-01)
$ws.Range("B320").Value = 6808405000000
$ws.Range("D320").Value = 354419833420.0937

# Row 327 (2023-03-01)
$ws.Range("B327").Value = 7965913000000
$ws.Range("D327").Value = 259054081300.813

# Row 329 (2023-05-01)
$ws.Range("B329").Value = 8140532000000
$ws.Range("D329").Value = 263874619124.7974

# --- New row 330 (2023-06-01) ---
# Copy the formatting of row 329 down to row 330 first (so the date cell
# keeps the same style, e.g. the date number format), then fill in values.
$ws.Range("A329:D329").Copy()
$ws.Range("A330:D330").PasteSpecial(-4122)

$ws.Range("A330").Value = 45078
$ws.Range("B330").Value = 8248190000000
$ws.Range("C330").Value = 0.03241491085899514
$ws.Range("D330").Value = 267364343598.0551
